$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# L5 was empty; it now holds the "xsd:uri" type annotation (same style s="2" as before)
$ws.Range("L5").Value = "xsd:uri"

# L8 held the constant year URI; bump it to 2016 and align its style with L7 (s="2")
# by copying L7's formatting onto L8 before updating the text, so the now-unused
# style slot (s="5") disappears on save.
$ws.Range("L7").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("L8").Value = "<http://reference.data.gov.uk/id/year/2016>"

# Row 11 ("Clase vivienda" duplicated under D11/G11) is removed entirely.
$ws.Rows("11:11").Delete()
